# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAllTypes = $wb.Worksheets.Item("全部类型")

# 展览 sheet (sheet1) updates
$wsExhibit.Range("F4").Value = 2948
$wsExhibit.Range("F6").Value = 107
$wsExhibit.Range("F8").Value = 1625
$wsExhibit.Range("F11").Value = 348
$wsExhibit.Range("F14").Value = 186
$wsExhibit.Range("F23").Value = 348
$wsExhibit.Range("F24").Value = 133
$wsExhibit.Range("F26").Value = 14
$wsExhibit.Range("F27").Value = 1963
$wsExhibit.Range("F29").Value = 450
$wsExhibit.Range("F30").Value = 10
$wsExhibit.Range("F31").Value = 157
$wsExhibit.Range("F33").Value = 218
$wsExhibit.Range("F34").Value = 327
$wsExhibit.Range("F36").Value = 480
$wsExhibit.Range("F37").Value = 4

# 全部类型 sheet (sheet4) updates
$wsAllTypes.Range("F4").Value = 2948
$wsAllTypes.Range("F6").Value = 107
$wsAllTypes.Range("F8").Value = 1625
$wsAllTypes.Range("F11").Value = 348
$wsAllTypes.Range("F14").Value = 186
$wsAllTypes.Range("F23").Value = 349
$wsAllTypes.Range("F24").Value = 133
$wsAllTypes.Range("F26").Value = 14
$wsAllTypes.Range("F27").Value = 1963
$wsAllTypes.Range("F29").Value = 450
$wsAllTypes.Range("F30").Value = 10
$wsAllTypes.Range("F31").Value = 157
$wsAllTypes.Range("F33").Value = 218
$wsAllTypes.Range("F34").Value = 327
$wsAllTypes.Range("F36").Value = 480
$wsAllTypes.Range("F37").Value = 4
